# Horarios actualizados Línea 141 - 387
# Updates the LP1912 / LP1912-215 / 6203-6173 scrape sheets with a fresh
# scrape pass ("Última actualización" 12:54:06 -> 13:23:09) that:
#  - re-shuffles a handful of existing "Hora_Scrap" ties (same arrival time,
#    different scrape timestamp) on sheet LP1912, and
#  - appends newly scraped arrivals (sorted back into place by Hora_Llegada)
#    on sheet LP1912 (rows 226-251) and sheet LP1912-215 (row 31).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# Sheet 1: LP1912
# ---------------------------------------------------------------------
$ws1 = $wb.Worksheets.Item("LP1912")

$ws1.Cells.Item(2, 1).Value = "Última actualización: 13:23:09"
$ws1.Cells.Item(3, 1).Value = "Total filas: 246"

# Each entry: row, Hora_Scrap, Hora_Llegada, Linea, Minutos, Parada
$sheet1Rows = @(
    @(35, "05:49:10", "07:05", "23_HERNANDEZ", 76, "LP1912"),
    @(36, "05:19:24", "07:05", "15_ABASTO", 106, "LP1912"),
    @(98, "08:48:09", "09:13", "16_SANTA ANA", 25, "LP1912"),
    @(99, "08:36:20", "09:13", "10_OLMOS", 37, "LP1912"),
    @(128, "09:25:56", "10:22", "16_SANTA ANA", 57, "LP1912"),
    @(129, "08:36:20", "10:22", "17_ROMERO", 106, "LP1912"),
    @(142, "10:12:35", "10:53", "10_OLMOS", 41, "LP1912"),
    @(143, "10:52:48", "10:53", "16_SANTA ANA", 1, "LP1912"),
    @(144, "09:25:56", "10:53", "27_EL RETIRO", 88, "LP1912"),
    @(175, "10:52:48", "11:58", "225_GOMEZ", 66, "LP1912"),
    @(176, "10:12:35", "11:58", "16_P MOR-167 Y 521", 106, "LP1912"),
    @(190, "11:17:08", "12:09", "10_OLMOS", 52, "LP1912"),
    @(191, "10:52:48", "12:09", "84_COLONIA URQUIZA-ESC 49", 77, "LP1912"),
    @(226, "13:23:09", "13:24", "16_SANTA ANA", 1, "LP1912"),
    @(227, "12:54:06", "13:26", "15_ABASTO", 32, "LP1912"),
    @(228, "11:45:01", "13:27", "14_ABASTO", 102, "LP1912"),
    @(229, "11:59:06", "13:31", "17_ROMERO", 92, "LP1912"),
    @(230, "12:27:08", "13:32", "10_OLMOS", 65, "LP1912"),
    @(231, "12:54:06", "13:35", "23_HERNANDEZ", 41, "LP1912"),
    @(232, "13:23:09", "13:36", "15_ABASTO", 13, "LP1912"),
    @(233, "12:27:08", "13:37", "23_HERNANDEZ", 70, "LP1912"),
    @(234, "12:27:08", "13:46", "17_ROMERO", 79, "LP1912"),
    @(235, "12:54:06", "13:47", "16_SANTA ANA", 53, "LP1912"),
    @(236, "12:54:06", "13:51", "11_ETCHEVERRY", 57, "LP1912"),
    @(237, "11:59:06", "13:51", "215A_EL PATO", 112, "LP1912"),
    @(238, "11:59:06", "13:56", "225_GOMEZ", 117, "LP1912"),
    @(239, "11:59:06", "13:57", "16_P MOR-167 Y 521", 118, "LP1912"),
    @(240, "12:27:08", "14:04", "17_ROMERO", 97, "LP1912"),
    @(241, "13:23:09", "14:05", "23_HERNANDEZ", 42, "LP1912"),
    @(242, "12:27:08", "14:17", "27_EL RETIRO", 110, "LP1912"),
    @(243, "12:27:08", "14:20", "215C_EL PATO", 113, "LP1912"),
    @(244, "12:54:06", "14:21", "26_HERNANDEZ", 87, "LP1912"),
    @(245, "12:54:06", "14:39", "14_ABASTO", 105, "LP1912"),
    @(246, "13:23:09", "14:57", "16_P MOR-SANTA ANA", 94, "LP1912"),
    @(247, "13:23:09", "14:58", "215B_EL PATO", 95, "LP1912"),
    @(248, "13:23:09", "15:00", "81_EL PELIGRO", 97, "LP1912"),
    @(249, "13:23:09", "15:05", "10_OLMOS", 102, "LP1912"),
    @(250, "13:23:09", "15:14", "11_ETCHEVERRY", 111, "LP1912"),
    @(251, "13:23:09", "15:21", "26_HERNANDEZ", 118, "LP1912")
)

foreach ($row in $sheet1Rows) {
    $r = $row[0]
    $ws1.Cells.Item($r, 1).Value = $row[1]
    $ws1.Cells.Item($r, 2).Value = $row[2]
    $ws1.Cells.Item($r, 3).Value = $row[3]
    $ws1.Cells.Item($r, 4).Value = $row[4]
    $ws1.Cells.Item($r, 5).Value = $row[5]
}

# ---------------------------------------------------------------------
# Sheet 2: LP1912-215
# ---------------------------------------------------------------------
$ws2 = $wb.Worksheets.Item("LP1912-215")

$ws2.Cells.Item(2, 1).Value = "Última actualización: 13:23:09"
$ws2.Cells.Item(3, 1).Value = "Total filas: 26"

$ws2.Cells.Item(31, 1).Value = "13:23:09"
$ws2.Cells.Item(31, 2).Value = "14:58"
$ws2.Cells.Item(31, 3).Value = "215B_EL PATO"
$ws2.Cells.Item(31, 4).Value = 95
$ws2.Cells.Item(31, 5).Value = "LP1912"

# ---------------------------------------------------------------------
# Sheet 3: 6203-6173
# ---------------------------------------------------------------------
$ws3 = $wb.Worksheets.Item("6203-6173")

$ws3.Cells.Item(2, 1).Value = "Última actualización: 13:23:09"
